# NMDC-EDGE-Metagenomics-ReadsQC bulk-submission template update
# - Reorders columns D:G (content + widths) so that:
#     old G (Sequencing Platform)                         -> new D
#     old D (Interleaved/Single-end Illumina/PacBio FASTQ) -> new E
#     old E (Illumina Paired-end R1 FASTQ)                 -> new F
#     old F (Illumina Paired-end R2 FASTQ)                 -> new G
# - Extends the hidden _xlnm._FilterDatabase defined name from A1:F1 to A1:G1
# - Rebuilds the data validations so sqref ranges line up with the new
#   column layout, merges the stray "F2 F3" validation back into the
#   Illumina R2 FASTQ validation (now G2:G99), extends the Data Source
#   validation to also cover column D (C2:D108) and adds error
#   alert title/message text to the Sequencing Platform and Data Source
#   dropdown validations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move column G in front of column D (shifts D,E,F -> E,F,G) ---
# This carries header cell values AND column widths along with it, matching
# Excel's native "drag column to new position" behaviour.
$ws.Columns.Item(7).Cut() | Out-Null
$ws.Columns.Item(4).Insert(-4161) | Out-Null

# --- 2. Widen the hidden AutoFilter defined name to include column G ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$1"
    }
}

# --- 3. Rebuild data validations to match the new layout ---
# Clear every existing validation on the sheet (their sqref ranges still
# refer to the pre-move column layout and need to be redone anyway).
$ws.Cells.Validation.Delete() | Out-Null

$illuminaR1Prompt = "Omit this column if input is Single fastq_x000a__x000a_Enter file name if Data Source is Uploaded File or Retrieved SRA Data_x000a__x000a_Enter file url if Data Source is HTTP(s) URL_x000a__x000a_Separate multiple inputs with commas"
$illuminaR2Prompt = $illuminaR1Prompt
$singleFastqPrompt = "Omit this column if input is Paired-end fastq_x000a__x000a_Enter file name if Data Source is Uploaded File or Retrieved SRA Data_x000a__x000a_Enter file url if Data Source is HTTP(s) URL_x000a__x000a_Separate multiple inputs with commas"

# Illumina R1 FASTQ -> F2:F99
$r = $ws.Range("F2:F99")
$r.Validation.Add(0, 1, 1, $null) | Out-Null
$r.Validation.InputTitle = "Illumina R1 FASTQ"
$r.Validation.InputMessage = $illuminaR1Prompt

# Project/Run Name -> A2:A99 (textLength 3-30, unchanged)
$r = $ws.Range("A2:A99")
$r.Validation.Add(6, 1, 1, 3, 30) | Out-Null
$r.Validation.ErrorTitle = "Project/Run Name"
$r.Validation.ErrorMessage = "Invalid Input"
$r.Validation.InputTitle = "Project/Run Name"
$r.Validation.InputMessage = "Required. At least 3 but less than 30 characters. Only alphabets, numbers, dashs, dot and underscore are allowed."

# Description -> B2:B99 (unchanged)
$r = $ws.Range("B2:B99")
$r.Validation.Add(0, 1, 1, $null) | Out-Null
$r.Validation.InputTitle = "Description"
$r.Validation.InputMessage = "Optional"

# Illumina R2 FASTQ -> G2:G99 (now one clean range, folding in the old F2/F3 remnant)
$r = $ws.Range("G2:G99")
$r.Validation.Add(0, 1, 1, $null) | Out-Null
$r.Validation.InputTitle = "Illumina R2 FASTQ"
$r.Validation.InputMessage = $illuminaR2Prompt

# Single Illumina/PacBio FASTQ -> E2:E99
$r = $ws.Range("E2:E99")
$r.Validation.Add(0, 1, 1, $null) | Out-Null
$r.Validation.InputTitle = "Single Illumina/PacBio FASTQ"
$r.Validation.InputMessage = $singleFastqPrompt

# Sequencing Platform -> D2:D99 (list, now also has error alert text)
$r = $ws.Range("D2:D99")
$r.Validation.Add(3, 1, 1, '"Illumina, PacBio"') | Out-Null
$r.Validation.ErrorTitle = "Sequencing Platform"
$r.Validation.ErrorMessage = "Select from dropdown list"
$r.Validation.InputTitle = "Sequencing Platform"
$r.Validation.InputMessage = "Default: Illumina"

# Data Source -> C2:D108 (list, now spans C & D, with error alert text)
$r = $ws.Range("C2:D108")
$r.Validation.Add(3, 1, 1, '"Uploaded File, Retrieved SRA Data, HTTP(s) URL"') | Out-Null
$r.Validation.ErrorTitle = "Data Source"
$r.Validation.ErrorMessage = "Select from dropdown list"
$r.Validation.InputTitle = "Data Source"
$r.Validation.InputMessage = "Default: Uploaded File"

"done"
